$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "-"
$ws.Range("E2").Value = "['MEC-2B-Des. Maq. Cad._T1', 'MEC-2B-Des. Maq. Cad._T1']"
$ws.Range("F2").Value = "-"

$ws.Range("C3").Value = "[-, -, 'MEC-3B-Cont.Lóg.Prog CLP', -]"

$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "[-, -, 'MEC-3B-Cont.Lóg.Prog CLP', -]"

$ws.Range("C6").Value = "[-, -, 'MEC-3B-Cont.Lóg.Prog CLP', -]"

$ws.Range("B7").Value = "['MEC-2B-Des. Maq. Cad._T2', 'MEC-2B-Des. Maq. Cad._T2']"
$ws.Range("C7").Value = "[-, -, 'MEC-3B-Cont.Lóg.Prog CLP', -]"
$ws.Range("F7").Value = "-"

$ws.Range("C8").Value = "-"
$ws.Range("E8").Value = "['MEC-2B-Des. Maq. Cad._T2', 'MEC-2B-Des. Maq. Cad._T1']"
$ws.Range("F8").Value = "-"

$ws.Range("B10").Value = "[-, 'MEC-1A-Desenho tecnico mecanico']"

$ws.Range("C12").Value = "-"

$ws.Range("C14").Value = "-"

$ws.Range("B16").Value = "['MEC-1A-Desenho tecnico mecanico', 'MEC-1A-Desenho tecnico mecanico']"

$ws.Range("C18").Value = "[-, 'MEC-2NB-C.pneumática', -, -]"
$ws.Range("D18").Value = "['MEC-1NA-Desenho tecnico mecanico – T2', 'MEC-1NA-Desenho tecnico mecanico – T2']"

$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "[-, 'MEC-2NB-C.pneumática', -, -]"
$ws.Range("D19").Value = "['MEC-1NA-Desenho tecnico mecanico – T2', 'MEC-1NA-Desenho tecnico mecanico – T2']"
$ws.Range("E19").Value = "-"

$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "[-, 'MEC-2NB-C.pneumática', -, -]"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "-"

$ws.Range("C21").Value = "[-, Elcio Dec.-C.pneumática-2NB, -, -]"
$ws.Range("D21").Value = "-"

$wb.Save()
